$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 191-192. This pushes the previous rows 191..229
# down to 193..231 (matching the target dimension A1:R231), preserving
# their formatting (incl. the date-number style on column D).
$ws.Rows("191:192").Insert()

# --- New row 191 ---------------------------------------------------
# Same record as the (now shifted) former row 191 (now at row 193),
# except the date moves from 2021-09-29 to 2021-10-07.
$ws.Range("A191").Value = $ws.Range("A193").Value2
$ws.Range("B191").Value = $ws.Range("B193").Value2
$ws.Range("C191").Value = $ws.Range("C193").Value2
$ws.Range("D191").Value = 44476
$ws.Range("E191").Value = $ws.Range("E193").Value2
$ws.Range("F191").Value = $ws.Range("F193").Value2
$ws.Range("G191").Value = $ws.Range("G193").Value2
$ws.Range("H191").Value = $ws.Range("H193").Value2
$ws.Range("I191").Value = $ws.Range("I193").Value2
$ws.Range("J191").Value = $ws.Range("J193").Value2
$ws.Range("K191").Value = $ws.Range("K193").Value2
$ws.Range("L191").Value = $ws.Range("L193").Value2
$ws.Range("M191").Value = $ws.Range("M193").Value2
$ws.Range("N191").Value = $ws.Range("N193").Value2
$ws.Range("O191").Value = $ws.Range("O193").Value2
$ws.Range("P191").Value = $ws.Range("P193").Value2
$ws.Range("Q191").Value = $ws.Range("Q193").Value2
$ws.Range("R191").Value = $ws.Range("R193").Value2

# --- New row 192 ---------------------------------------------------
# Same record as the (now shifted) former row 192 (now at row 194),
# except date, volume, unit, origin and the two derived columns change.
$ws.Range("A192").Value = $ws.Range("A194").Value2
$ws.Range("B192").Value = $ws.Range("B194").Value2
$ws.Range("C192").Value = $ws.Range("C194").Value2
$ws.Range("D192").Value = 44476
$ws.Range("E192").Value = $ws.Range("E194").Value2
$ws.Range("F192").Value = $ws.Range("F194").Value2
$ws.Range("G192").Value = $ws.Range("G194").Value2
$ws.Range("H192").Value = $ws.Range("H194").Value2
$ws.Range("I192").Value = $ws.Range("I194").Value2
$ws.Range("J192").Value = 150
$ws.Range("K192").Value = $ws.Range("K194").Value2
$ws.Range("L192").Value = $ws.Range("L194").Value2
$ws.Range("M192").Value = $ws.Range("M194").Value2
$ws.Range("N192").Value = "`$/caja 80 unidades"
$ws.Range("O192").Value = "Región del Maule"
$ws.Range("P192").Value = 188
$ws.Range("Q192").Value = 80
$ws.Range("R192").Value = $ws.Range("R194").Value2
